# Apply a cyclic rotation of the weekly price records between rows 2, 3 and 4:
#   new row2 = old row4
#   new row3 = old row2
#   new row4 = old row3
# Only columns D, J, K, L, M, O, P change (per the diff); the rest stay put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("D", "J", "K", "L", "M", "O", "P")

# Capture current ("before") values for the rows involved.
$before = @{}
foreach ($r in 2, 3, 4) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range("$col$r").Value2
    }
    $before[$r] = $rowVals
}

# Mapping: after row2 <- before row4 ; after row3 <- before row2 ; after row4 <- before row3
$source = @{ 2 = 4; 3 = 2; 4 = 3 }

foreach ($r in 2, 3, 4) {
    $src = $source[$r]
    foreach ($col in $cols) {
        $ws.Range("$col$r").Value = $before[$src][$col]
    }
}
